$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new Key/Locators rows for the Reviews & Ratings test cases ---
$ws.Range("A25").Value = "Submit_Review_Xpath"
$ws.Range("B25").Value = "//*[@id='review-form']/fieldset/div[2]/div[2]/button"

$ws.Range("A26").Value = "Review_Title_Textbox_ErrorMsg_Xpath"
$ws.Range("B26").Value = "//*[@id='advice-required-entry-summary_field']"

$ws.Range("A27").Value = "Review_Textbox_ErrorMsg_Xpath"
$ws.Range("B27").Value = "//*[@id='advice-required-entry-review_field']"

# --- Tidy up A16's leftover duplicate "blue Courier New" font formatting by
#     re-applying its own format, which collapses it onto the equivalent,
#     already-used style instead of the redundant one. ---
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

# --- Move the active selection down below the newly entered rows ---
$ws.Range("A29").Select()
